$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4383.3335
$ws.Range("I76").Value = 3800
$ws.Range("J76").Value = 4500
$ws.Range("K76").Value = 3800
$ws.Range("L76").Value = 4500
$ws.Range("M76").Value = -3485
$ws.Range("N76").Value = -5130
$ws.Range("H79").Value = 4383.3335
$ws.Range("I79").Value = 3800
$ws.Range("J79").Value = 4500
$ws.Range("K79").Value = 3800
$ws.Range("L79").Value = 4500
$ws.Range("M79").Value = -2708
$ws.Range("N79").Value = -6684
$ws.Range("H113").Value = 2793.0952
$ws.Range("I113").Value = 2425.625
$ws.Range("J113").Value = 3019.2307
$ws.Range("K113").Value = 2425.625
$ws.Range("L113").Value = 3019.2307
$ws.Range("M113").Value = 828.375
$ws.Range("N113").Value = -9527.2307
$ws.Range("H137").Value = 1607.4054
$ws.Range("I137").Value = 1061.5555
$ws.Range("J137").Value = 3081.2
$ws.Range("K137").Value = 3184.6665
$ws.Range("L137").Value = 9243.599999999999
$ws.Range("M137").Value = -634.6664999999998
$ws.Range("N137").Value = -14343.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 777.625
$ws.Range("I2").Value = 760.1429000000001
$ws.Range("K2").Value = 760.1429000000001
$ws.Range("M2").Value = -647.1429000000001
$ws.Range("H45").Value = 2227.6667
$ws.Range("J45").Value = 1967
$ws.Range("L45").Value = 1967
$ws.Range("N45").Value = -2721
$ws.Range("H74").Value = 4809720
$ws.Range("I74").Value = 8622769
$ws.Range("J74").Value = 1962.5217
$ws.Range("K74").Value = 8622769
$ws.Range("L74").Value = 1962.5217
$ws.Range("M74").Value = -8621895
$ws.Range("N74").Value = -3710.5217
$ws.Range("H77").Value = 4809720
$ws.Range("I77").Value = 8622769
$ws.Range("J77").Value = 1962.5217
$ws.Range("K77").Value = 43113845
$ws.Range("L77").Value = 9812.6085
$ws.Range("M77").Value = -43109477
$ws.Range("N77").Value = -18548.6085
$ws.Range("H116").Value = 777.625
$ws.Range("I116").Value = 760.1429000000001
$ws.Range("K116").Value = 760.1429000000001
$ws.Range("M116").Value = 1533.8571
$ws.Range("H122").Value = 2464.9092
$ws.Range("I122").Value = 1825.5
$ws.Range("K122").Value = 5476.5
$ws.Range("M122").Value = -3026.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 777.625
$ws.Range("I3").Value = 760.1429000000001
$ws.Range("K3").Value = 760.1429000000001
$ws.Range("M3").Value = -646.1429000000001
$ws.Range("H80").Value = 45.5625
$ws.Range("J80").Value = 51.214287
$ws.Range("L80").Value = 51.214287
$ws.Range("N80").Value = -2047.214287
$ws.Range("H83").Value = 45.5625
$ws.Range("J83").Value = 51.214287
$ws.Range("L83").Value = 256.071435
$ws.Range("N83").Value = -10240.071435
$ws.Range("H94").Value = 40352.2
$ws.Range("I94").Value = 429.75
$ws.Range("K94").Value = 429.75
$ws.Range("M94").Value = 21.25
$ws.Range("H107").Value = 1672.8636
$ws.Range("I107").Value = 1300.6666
$ws.Range("K107").Value = 1300.6666
$ws.Range("M107").Value = 619.3334

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3093636.5
$ws.Range("I31").Value = 4903318
$ws.Range("J31").Value = 2096.875
$ws.Range("K31").Value = 4903318
$ws.Range("L31").Value = 2096.875
$ws.Range("M31").Value = -4903023
$ws.Range("N31").Value = -2686.875
$ws.Range("H34").Value = 3093636.5
$ws.Range("I34").Value = 4903318
$ws.Range("J34").Value = 2096.875
$ws.Range("K34").Value = 4903318
$ws.Range("L34").Value = 2096.875
$ws.Range("M34").Value = -4903116
$ws.Range("N34").Value = -2500.875
$ws.Range("H80").Value = 16055.556
$ws.Range("J80").Value = 16055.556
$ws.Range("L80").Value = 16055.556
$ws.Range("N80").Value = -18301.556
$ws.Range("H83").Value = 16055.556
$ws.Range("J83").Value = 16055.556
$ws.Range("L83").Value = 48166.66800000001
$ws.Range("N83").Value = -59398.66800000001
$ws.Range("H86").Value = 386913.16
$ws.Range("I86").Value = 626695.4
$ws.Range("J86").Value = 3261.6
$ws.Range("K86").Value = 626695.4
$ws.Range("L86").Value = 3261.6
$ws.Range("M86").Value = -625572.4
$ws.Range("N86").Value = -5507.6
$ws.Range("H88").Value = 32049
$ws.Range("J88").Value = 32049
$ws.Range("L88").Value = 32049
$ws.Range("N88").Value = -32861
$ws.Range("H89").Value = 386913.16
$ws.Range("I89").Value = 626695.4
$ws.Range("J89").Value = 3261.6
$ws.Range("K89").Value = 3133477
$ws.Range("L89").Value = 16308
$ws.Range("M89").Value = -3127861
$ws.Range("N89").Value = -27540
$ws.Range("H91").Value = 32049
$ws.Range("J91").Value = 32049
$ws.Range("L91").Value = 32049
$ws.Range("N91").Value = -34857
$ws.Range("H107").Value = 33335088
$ws.Range("I107").Value = 62502068
$ws.Range("J107").Value = 1394.7142
$ws.Range("K107").Value = 62502068
$ws.Range("L107").Value = 1394.7142
$ws.Range("M107").Value = -62500148
$ws.Range("N107").Value = -5234.7142
$ws.Range("H122").Value = 16667666
$ws.Range("I122").Value = 22223222
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 66669666
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -66667216
$ws.Range("N122").Value = -7900
$ws.Range("H132").Value = 1484.05
$ws.Range("I132").Value = 1190.2778
$ws.Range("J132").Value = 4128
$ws.Range("K132").Value = 3570.8334
$ws.Range("L132").Value = 12384
$ws.Range("M132").Value = -1040.8334
$ws.Range("N132").Value = -17444
$ws.Range("H134").Value = 5202.8
$ws.Range("I134").Value = 5992.1816
$ws.Range("J134").Value = 3032
$ws.Range("K134").Value = 17976.5448
$ws.Range("L134").Value = 9096
$ws.Range("M134").Value = -15441.5448
$ws.Range("N134").Value = -14166

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 243.13333
$ws.Range("I8").Value = 243.13333
$ws.Range("K8").Value = 729.39999
$ws.Range("M8").Value = -590.39999
$ws.Range("H113").Value = 1376.0769
$ws.Range("I113").Value = 612
$ws.Range("J113").Value = 2598.6
$ws.Range("K113").Value = 1836
$ws.Range("L113").Value = 7795.799999999999
$ws.Range("M113").Value = 334
$ws.Range("N113").Value = -12135.8
$ws.Range("H131").Value = 13889764
$ws.Range("J131").Value = 17544790
$ws.Range("L131").Value = 52634370
$ws.Range("N131").Value = -52644450

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 59800
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 59800
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 59800
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -61672
$ws.Range("H77").Value = 59800
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 59800
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 179400
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -188760
$ws.Range("H80").Value = 3497.0334
$ws.Range("I80").Value = 5867.5
$ws.Range("J80").Value = 2904.4167
$ws.Range("K80").Value = 5867.5
$ws.Range("L80").Value = 2904.4167
$ws.Range("M80").Value = -4869.5
$ws.Range("N80").Value = -4900.4167
$ws.Range("H83").Value = 3497.0334
$ws.Range("I83").Value = 5867.5
$ws.Range("J83").Value = 2904.4167
$ws.Range("K83").Value = 29337.5
$ws.Range("L83").Value = 14522.0835
$ws.Range("M83").Value = -24345.5
$ws.Range("N83").Value = -24506.0835
$ws.Range("H113").Value = 1439.6086
$ws.Range("I113").Value = 1012.75
$ws.Range("J113").Value = 1905.2727
$ws.Range("K113").Value = 1012.75
$ws.Range("L113").Value = 1905.2727
$ws.Range("M113").Value = 1157.25
$ws.Range("N113").Value = -6245.2727
$ws.Range("H122").Value = 3934.5833
$ws.Range("I122").Value = 3506.6875
$ws.Range("K122").Value = 10520.0625
$ws.Range("M122").Value = -8070.0625
$ws.Range("H126").Value = 1748.8
$ws.Range("I126").Value = 1568.7142
$ws.Range("K126").Value = 4706.142599999999
$ws.Range("M126").Value = -2236.142599999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 750.0625
$ws.Range("I22").Value = 407.14285
$ws.Range("K22").Value = 407.14285
$ws.Range("M22").Value = -112.14285
$ws.Range("H27").Value = 750.0625
$ws.Range("I27").Value = 407.14285
$ws.Range("K27").Value = 407.14285
$ws.Range("M27").Value = -300.14285
$ws.Range("H40").Value = 6080.8
$ws.Range("I40").Value = 2301.3333
$ws.Range("J40").Value = 11750
$ws.Range("K40").Value = 2301.3333
$ws.Range("L40").Value = 11750
$ws.Range("M40").Value = -2165.3333
$ws.Range("N40").Value = -12022
$ws.Range("H55").Value = 496.83334
$ws.Range("I55").Value = 527
$ws.Range("K55").Value = 527
$ws.Range("M55").Value = -354
$ws.Range("H122").Value = 3975.5715
$ws.Range("I122").Value = 4999.75
$ws.Range("J122").Value = 2610
$ws.Range("K122").Value = 14999.25
$ws.Range("L122").Value = 7830
$ws.Range("M122").Value = -12549.25
$ws.Range("N122").Value = -12730

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5233.273
$ws.Range("I62").Value = 5233.3335
$ws.Range("J62").Value = 5233.25
$ws.Range("K62").Value = 5233.3335
$ws.Range("L62").Value = 5233.25
$ws.Range("M62").Value = -4609.3335
$ws.Range("N62").Value = -6481.25
$ws.Range("H65").Value = 5233.273
$ws.Range("I65").Value = 5233.3335
$ws.Range("J65").Value = 5233.25
$ws.Range("K65").Value = 26166.6675
$ws.Range("L65").Value = 26166.25
$ws.Range("M65").Value = -23046.6675
$ws.Range("N65").Value = -32406.25
$ws.Range("H81").Value = 100003110
$ws.Range("I81").Value = 250003460
$ws.Range("K81").Value = 500006920
$ws.Range("M81").Value = -500005859
$ws.Range("H84").Value = 100003110
$ws.Range("I84").Value = 250003460
$ws.Range("K84").Value = 2500034600
$ws.Range("M84").Value = -2500029296
$ws.Range("H122").Value = 76924800
$ws.Range("I122").Value = 83334950
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 250004850
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -250002400
$ws.Range("N122").Value = -13900
$ws.Range("H126").Value = 4740.4
$ws.Range("I126").Value = 5613
$ws.Range("K126").Value = 16839
$ws.Range("M126").Value = -14369
